$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.007798538995164189
$ws.Range("C2").Value = 0.8629205835071485
$ws.Range("D2").Value = 1.444608608967088
$ws.Range("E2").Value = 1.20191871978395
$ws.Range("F2").Value = 1.215474567149581

# Row 3 (Q0)
$ws.Range("B3").Value = 0.09720554110293481
$ws.Range("C3").Value = 1.198028828389808
$ws.Range("D3").Value = 3.637881730776367
$ws.Range("E3").Value = 1.907323184669123
$ws.Range("F3").Value = 1.911493256486692
$ws.Range("G3").Value = 144

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1802234366250808
$ws.Range("C4").Value = 1.394212243063904
$ws.Range("D4").Value = 8.046150714213496
$ws.Range("E4").Value = 2.836573763224481
$ws.Range("F4").Value = 2.851282233543963
$ws.Range("G4").Value = 70
